$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name corrections (A column) & updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 12:46"
$ws.Range("A174").Value = "Namibia"
$ws.Range("A175").Value = "Bermudas"
$ws.Range("A176").Value = "Camboya"
$ws.Range("A177").Value = "Brunei"
$ws.Range("A195").Value = "Islas Turcas y Caicos"
$ws.Range("A196").Value = "Lesoto"
$ws.Range("A197").Value = "Belice"
$ws.Range("A198").Value = "Timor Oriental"
$ws.Range("A199").Value = "Curazao"
$ws.Range("A200").Value = "Granada"
$ws.Range("A201").Value = "Nueva Caledonia"
$ws.Range("A202").Value = "Seychelles"
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"
$ws.Range("A206").Value = "Dominica"
$ws.Range("A207").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Updated statistics (B..H columns) ---
$ws.Range("B4").Value = 2596771
$ws.Range("C4").Value = 234
$ws.Range("E4").Value = 1387125
$ws.Range("B7").Value = 530924
$ws.Range("C7").Value = 1347
$ws.Range("D7").Value = 311001
$ws.Range("E7").Value = 203799
$ws.Range("G7").Value = 21
$ws.Range("H7").Value = 16124
$ws.Range("B13").Value = 222669
$ws.Range("C13").Value = 2489
$ws.Range("D13").Value = 183310
$ws.Range("E13").Value = 28851
$ws.Range("G13").Value = 144
$ws.Range("H13").Value = 10508
$ws.Range("B43").Value = 31617
$ws.Range("C43").Value = 62
$ws.Range("E43").Value = 555
$ws.Range("B48").Value = 26313
$ws.Range("C48").Value = 291
$ws.Range("D48").Value = 18814
$ws.Range("E48").Value = 5887
$ws.Range("G48").Value = 23
$ws.Range("H48").Value = 1612
$ws.Range("B69").Value = 11306
$ws.Range("C69").Value = 8
$ws.Range("E69").Value = 3277
$ws.Range("H69").Value = 347
$ws.Range("B78").Value = 6586
$ws.Range("C78").Value = 127
$ws.Range("D78").Value = 4291
$ws.Range("E78").Value = 2190
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 105
$ws.Range("B103").Value = 2402
$ws.Range("C103").Value = 72
$ws.Range("D103").Value = 1384
$ws.Range("E103").Value = 963
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 55
$ws.Range("B124").Value = 1200
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 1104
$ws.Range("E124").Value = 89
$ws.Range("B126").Value = 1149
$ws.Range("C126").Value = 25
$ws.Range("D126").Value = 306
$ws.Range("E126").Value = 827
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 16
$ws.Range("B132").Value = 1068
$ws.Range("C132").Value = 6
$ws.Range("D132").Value = 937
$ws.Range("E132").Value = 64
$ws.Range("B174").Value = 150
$ws.Range("C174").Value = 14
$ws.Range("D174").Value = 22
$ws.Range("E174").Value = 128
$ws.Range("H174").Value = 0
$ws.Range("B175").Value = 146
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 134
$ws.Range("E175").Value = 3
$ws.Range("H175").Value = 9
$ws.Range("C176").Value = 2
$ws.Range("D176").Value = 129
$ws.Range("E176").Value = 12
$ws.Range("H176").Value = 0
$ws.Range("B177").Value = 141
$ws.Range("D177").Value = 138
$ws.Range("E177").Value = 0
$ws.Range("H177").Value = 3
$ws.Range("B195").Value = 28
$ws.Range("C195").Value = 12
$ws.Range("D195").Value = 11
$ws.Range("E195").Value = 16
$ws.Range("H195").Value = 1
$ws.Range("D196").Value = 4
$ws.Range("E196").Value = 20
$ws.Range("H196").Value = 0
$ws.Range("D197").Value = 18
$ws.Range("E197").Value = 4
$ws.Range("H197").Value = 2
$ws.Range("B198").Value = 24
$ws.Range("D198").Value = 24
$ws.Range("E198").Value = 0
$ws.Range("H198").Value = 0
$ws.Range("D199").Value = 19
$ws.Range("E199").Value = 3
$ws.Range("H199").Value = 1
$ws.Range("B200").Value = 23
$ws.Range("D200").Value = 23
$ws.Range("B201").Value = 21
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 21
$ws.Range("E201").Value = 0
$ws.Range("B202").Value = 20
$ws.Range("C202").Value = 9
$ws.Range("D202").Value = 11
$ws.Range("E202").Value = 9
$ws.Range("B204").Value = 19
$ws.Range("D204").Value = 19
$ws.Range("B206").Value = 18
$ws.Range("D206").Value = 18
$ws.Range("E206").Value = 0
$ws.Range("B207").Value = 17
$ws.Range("D207").Value = 0
$ws.Range("E207").Value = 17
$ws.Range("H207").Value = 0
